$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 390.19232
$ws.Range("I2").Value = 363.8889
$ws.Range("J2").Value = 449.375
$ws.Range("K2").Value = 363.8889
$ws.Range("L2").Value = 449.375
$ws.Range("M2").Value = -250.8889
$ws.Range("N2").Value = -675.375
$ws.Range("H9").Value = 16667436
$ws.Range("I9").Value = 37037172
$ws.Range("J9").Value = 1290.2727
$ws.Range("K9").Value = 37037172
$ws.Range("L9").Value = 1290.2727
$ws.Range("M9").Value = -37037003
$ws.Range("N9").Value = -1628.2727
$ws.Range("H15").Value = 796.5574
$ws.Range("I15").Value = 796.5574
$ws.Range("K15").Value = 2389.6722
$ws.Range("M15").Value = -2220.6722
$ws.Range("H17").Value = 376.53522
$ws.Range("J17").Value = 380.2029
$ws.Range("L17").Value = 1140.6087
$ws.Range("N17").Value = -1476.6087
$ws.Range("H64").Value = 5119.148
$ws.Range("I64").Value = 4653.722
$ws.Range("J64").Value = 6050
$ws.Range("K64").Value = 4653.722
$ws.Range("L64").Value = 6050
$ws.Range("M64").Value = -4405.722
$ws.Range("N64").Value = -6546
$ws.Range("H67").Value = 5119.148
$ws.Range("I67").Value = 4653.722
$ws.Range("J67").Value = 6050
$ws.Range("K67").Value = 4653.722
$ws.Range("L67").Value = 6050
$ws.Range("M67").Value = -3795.722
$ws.Range("N67").Value = -7766
$ws.Range("H106").Value = 2356.5833
$ws.Range("I106").Value = 2475.4443
$ws.Range("K106").Value = 2475.4443
$ws.Range("M106").Value = -1844.4443
$ws.Range("H119").Value = 2070
$ws.Range("J119").Value = 2070
$ws.Range("L119").Value = 6210
$ws.Range("N119").Value = -15886
$ws.Range("H134").Value = 95139.5
$ws.Range("J134").Value = 95139.5
$ws.Range("L134").Value = 95139.5
$ws.Range("N134").Value = -105279.5
$ws.Range("H135").Value = 1816.0625
$ws.Range("I135").Value = 1816.0625
$ws.Range("K135").Value = 16344.5625
$ws.Range("M135").Value = -13809.5625
$ws.Range("H136").Value = 77977.336
$ws.Range("J136").Value = 77977.336
$ws.Range("L136").Value = 77977.336
$ws.Range("N136").Value = -88177.336
$ws.Range("H137").Value = 384454.7
$ws.Range("I137").Value = 2155.6428
$ws.Range("J137").Value = 607462.4399999999
$ws.Range("K137").Value = 6466.928400000001
$ws.Range("L137").Value = 1822387.32
$ws.Range("M137").Value = -3916.928400000001
$ws.Range("N137").Value = -1827487.32
$ws.Range("H138").Value = 42277.6
$ws.Range("I138").Value = 68285.60000000001
$ws.Range("J138").Value = 3265.6
$ws.Range("K138").Value = 204856.8
$ws.Range("L138").Value = 9796.799999999999
$ws.Range("M138").Value = -199716.8
$ws.Range("N138").Value = -20076.8
$ws.Range("H140").Value = 83270.28999999999
$ws.Range("J140").Value = 92030.5
$ws.Range("L140").Value = 92030.5
$ws.Range("N140").Value = -102390.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2026.5
$ws.Range("I5").Value = 252
$ws.Range("K5").Value = 252
$ws.Range("M5").Value = -140
$ws.Range("H7").Value = 56869.1
$ws.Range("J7").Value = 58115.168
$ws.Range("L7").Value = 58115.168
$ws.Range("N7").Value = -58343.168
$ws.Range("H63").Value = 8119.2856
$ws.Range("I63").Value = 4750
$ws.Range("J63").Value = 9467
$ws.Range("K63").Value = 4750
$ws.Range("L63").Value = 9467
$ws.Range("M63").Value = -4064
$ws.Range("N63").Value = -10839
$ws.Range("H66").Value = 8119.2856
$ws.Range("I66").Value = 4750
$ws.Range("J66").Value = 9467
$ws.Range("K66").Value = 23750
$ws.Range("L66").Value = 47335
$ws.Range("M66").Value = -20318
$ws.Range("N66").Value = -54199
$ws.Range("H88").Value = 1601.6923
$ws.Range("I88").Value = 902.8
$ws.Range("J88").Value = 2038.5
$ws.Range("K88").Value = 902.8
$ws.Range("L88").Value = 2038.5
$ws.Range("M88").Value = -496.8
$ws.Range("N88").Value = -2850.5
$ws.Range("H91").Value = 1601.6923
$ws.Range("I91").Value = 902.8
$ws.Range("J91").Value = 2038.5
$ws.Range("K91").Value = 902.8
$ws.Range("L91").Value = 2038.5
$ws.Range("M91").Value = 501.2
$ws.Range("N91").Value = -4846.5
$ws.Range("H104").Value = 28515.334
$ws.Range("J104").Value = 28515.334
$ws.Range("L104").Value = 28515.334
$ws.Range("N104").Value = -35503.334
$ws.Range("H107").Value = 55191.6
$ws.Range("J107").Value = 55191.6
$ws.Range("L107").Value = 55191.6
$ws.Range("N107").Value = -62871.6
$ws.Range("H118").Value = 52612.5
$ws.Range("J118").Value = 52612.5
$ws.Range("L118").Value = 52612.5
$ws.Range("N118").Value = -55926.5
$ws.Range("H121").Value = 69684.09
$ws.Range("J121").Value = 75052.5
$ws.Range("L121").Value = 75052.5
$ws.Range("N121").Value = -78546.5
$ws.Range("H122").Value = 1214.1818
$ws.Range("I122").Value = 989.12
$ws.Range("J122").Value = 1917.5
$ws.Range("K122").Value = 2967.36
$ws.Range("L122").Value = 5752.5
$ws.Range("M122").Value = -517.3600000000001
$ws.Range("N122").Value = -10652.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 27491.908
$ws.Range("J2").Value = 27491.908
$ws.Range("L2").Value = 27491.908
$ws.Range("N2").Value = -27717.908
$ws.Range("H4").Value = 2026.5
$ws.Range("I4").Value = 252
$ws.Range("K4").Value = 252
$ws.Range("M4").Value = -137
$ws.Range("H13").Value = 53589.6
$ws.Range("J13").Value = 53589.6
$ws.Range("L13").Value = 53589.6
$ws.Range("N13").Value = -53925.6
$ws.Range("H50").Value = 53131.4
$ws.Range("J50").Value = 53131.4
$ws.Range("L50").Value = 53131.4
$ws.Range("N50").Value = -54279.4
$ws.Range("H51").Value = 40182.4
$ws.Range("J51").Value = 40182.4
$ws.Range("L51").Value = 40182.4
$ws.Range("N51").Value = -41164.4
$ws.Range("H53").Value = 33175.8
$ws.Range("J53").Value = 33175.8
$ws.Range("L53").Value = 33175.8
$ws.Range("N53").Value = -34323.8
$ws.Range("H55").Value = 28581.166
$ws.Range("J55").Value = 28581.166
$ws.Range("L55").Value = 28581.166
$ws.Range("N55").Value = -29127.166
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H86").Value = 3362.889
$ws.Range("I86").Value = 2389.9
$ws.Range("K86").Value = 2389.9
$ws.Range("M86").Value = -1266.9
$ws.Range("H88").Value = 44999.5
$ws.Range("J88").Value = 44999.5
$ws.Range("L88").Value = 44999.5
$ws.Range("N88").Value = -45811.5
$ws.Range("H89").Value = 3362.889
$ws.Range("I89").Value = 2389.9
$ws.Range("K89").Value = 11949.5
$ws.Range("M89").Value = -6333.5
$ws.Range("H91").Value = 44999.5
$ws.Range("J91").Value = 44999.5
$ws.Range("L91").Value = 44999.5
$ws.Range("N91").Value = -47807.5
$ws.Range("H114").Value = 99919.664
$ws.Range("J114").Value = 99919.664
$ws.Range("L114").Value = 99919.664
$ws.Range("N114").Value = -108597.664
$ws.Range("H116").Value = 99990
$ws.Range("J116").Value = 99990
$ws.Range("L116").Value = 99990
$ws.Range("N116").Value = -109168
$ws.Range("H117").Value = 98139.664
$ws.Range("J117").Value = 98139.664
$ws.Range("L117").Value = 98139.664
$ws.Range("N117").Value = -107317.664
$ws.Range("H118").Value = 72046.57000000001
$ws.Range("J118").Value = 74796.664
$ws.Range("L118").Value = 74796.664
$ws.Range("N118").Value = -78110.664
$ws.Range("H119").Value = 60996
$ws.Range("J119").Value = 60996
$ws.Range("L119").Value = 60996
$ws.Range("N119").Value = -70672

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 79995
$ws.Range("J9").Value = 79995
$ws.Range("L9").Value = 79995
$ws.Range("N9").Value = -80331
$ws.Range("H31").Value = 3906.7778
$ws.Range("I31").Value = 2037.3334
$ws.Range("K31").Value = 2037.3334
$ws.Range("M31").Value = -1742.3334
$ws.Range("H34").Value = 3906.7778
$ws.Range("I34").Value = 2037.3334
$ws.Range("K34").Value = 2037.3334
$ws.Range("M34").Value = -1835.3334
$ws.Range("H41").Value = 14200
$ws.Range("J41").Value = 15250
$ws.Range("L41").Value = 15250
$ws.Range("N41").Value = -16106
$ws.Range("H58").Value = 3447.738
$ws.Range("J58").Value = 3298.8572
$ws.Range("L58").Value = 3298.8572
$ws.Range("N58").Value = -3704.8572
$ws.Range("H60").Value = 31666.666
$ws.Range("I60").Value = 17500
$ws.Range("J60").Value = 60000
$ws.Range("K60").Value = 17500
$ws.Range("L60").Value = 60000
$ws.Range("M60").Value = -16989
$ws.Range("N60").Value = -61022
$ws.Range("H62").Value = 2552.3333
$ws.Range("I62").Value = 2269.6667
$ws.Range("K62").Value = 2269.6667
$ws.Range("M62").Value = -1645.6667
$ws.Range("H65").Value = 2552.3333
$ws.Range("I65").Value = 2269.6667
$ws.Range("K65").Value = 11348.3335
$ws.Range("M65").Value = -8228.333500000001
$ws.Range("H105").Value = 3958.4119
$ws.Range("I105").Value = 2499.4
$ws.Range("K105").Value = 2499.4
$ws.Range("M105").Value = -752.4000000000001
$ws.Range("H117").Value = 31249.5
$ws.Range("J117").Value = 31249.5
$ws.Range("L117").Value = 31249.5
$ws.Range("N117").Value = -40427.5
$ws.Range("H122").Value = 2640.2068
$ws.Range("I122").Value = 2061.0833
$ws.Range("J122").Value = 5420
$ws.Range("K122").Value = 6183.249899999999
$ws.Range("L122").Value = 16260
$ws.Range("M122").Value = -3733.249899999999
$ws.Range("N122").Value = -21160
$ws.Range("H132").Value = 2094948.9
$ws.Range("I132").Value = 1901750.4
$ws.Range("J132").Value = 3254140
$ws.Range("K132").Value = 5705251.199999999
$ws.Range("L132").Value = 9762420
$ws.Range("M132").Value = -5702721.199999999
$ws.Range("N132").Value = -9767480
$ws.Range("H133").Value = 102747.5
$ws.Range("I133").Value = 85000
$ws.Range("J133").Value = 108663.336
$ws.Range("K133").Value = 85000
$ws.Range("L133").Value = 108663.336
$ws.Range("M133").Value = -82470
$ws.Range("N133").Value = -113723.336
$ws.Range("H134").Value = 3762456.8
$ws.Range("I134").Value = 5497572
$ws.Range("J134").Value = 3040.5
$ws.Range("K134").Value = 16492716
$ws.Range("L134").Value = 9121.5
$ws.Range("M134").Value = -16490181
$ws.Range("N134").Value = -14191.5
$ws.Range("H136").Value = 3447.738
$ws.Range("J136").Value = 3298.8572
$ws.Range("L136").Value = 9896.571599999999
$ws.Range("N136").Value = -14996.5716

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 6047.4
$ws.Range("J32").Value = 2670.3572
$ws.Range("L32").Value = 8011.071599999999
$ws.Range("N32").Value = -8577.071599999999
$ws.Range("H38").Value = 187.60715
$ws.Range("I38").Value = 9.608696
$ws.Range("J38").Value = 1006.4
$ws.Range("K38").Value = 28.826088
$ws.Range("L38").Value = 3019.2
$ws.Range("M38").Value = 318.173912
$ws.Range("N38").Value = -3713.2
$ws.Range("H113").Value = 6078250.5
$ws.Range("I113").Value = 1164.25
$ws.Range("K113").Value = 3492.75
$ws.Range("M113").Value = -1322.75
$ws.Range("H121").Value = 771995.75
$ws.Range("I121").Value = 1741
$ws.Range("K121").Value = 5223
$ws.Range("M121").Value = -3913
$ws.Range("H131").Value = 1780.8889
$ws.Range("J131").Value = 2171.3333
$ws.Range("L131").Value = 6513.999899999999
$ws.Range("N131").Value = -16593.9999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 565.4815
$ws.Range("I107").Value = 561.6957
$ws.Range("J107").Value = 587.25
$ws.Range("K107").Value = 561.6957
$ws.Range("L107").Value = 587.25
$ws.Range("M107").Value = 1358.3043
$ws.Range("N107").Value = -4427.25
$ws.Range("H109").Value = 23388.818
$ws.Range("J109").Value = 23388.818
$ws.Range("L109").Value = 23388.818
$ws.Range("N109").Value = -25468.818
$ws.Range("H122").Value = 4971.96
$ws.Range("I122").Value = 4565
$ws.Range("K122").Value = 13695
$ws.Range("M122").Value = -11245
$ws.Range("H132").Value = 7765.3076
$ws.Range("I132").Value = 2605.4443
$ws.Range("K132").Value = 7816.3329
$ws.Range("M132").Value = -5286.3329

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12849.272
$ws.Range("I7").Value = 36292.668
$ws.Range("J7").Value = 4058
$ws.Range("K7").Value = 36292.668
$ws.Range("L7").Value = 4058
$ws.Range("M7").Value = -36180.668
$ws.Range("N7").Value = -4282
$ws.Range("H61").Value = 3143.6155
$ws.Range("I61").Value = 3424.4546
$ws.Range("K61").Value = 3424.4546
$ws.Range("M61").Value = -3222.4546
$ws.Range("H68").Value = 3770.2
$ws.Range("I68").Value = 3833.1667
$ws.Range("K68").Value = 3833.1667
$ws.Range("M68").Value = -3084.1667
$ws.Range("H71").Value = 3770.2
$ws.Range("I71").Value = 3833.1667
$ws.Range("K71").Value = 19165.8335
$ws.Range("M71").Value = -15421.8335
$ws.Range("H113").Value = 3143.6155
$ws.Range("I113").Value = 3424.4546
$ws.Range("K113").Value = 3424.4546
$ws.Range("M113").Value = -1254.4546
$ws.Range("H117").Value = 57595.5
$ws.Range("J117").Value = 57595.5
$ws.Range("L117").Value = 57595.5
$ws.Range("N117").Value = -66773.5
$ws.Range("H121").Value = 48208.418
$ws.Range("J121").Value = 48208.418
$ws.Range("L121").Value = 48208.418
$ws.Range("N121").Value = -51702.418
$ws.Range("H126").Value = 12849.272
$ws.Range("I126").Value = 36292.668
$ws.Range("J126").Value = 4058
$ws.Range("K126").Value = 108878.004
$ws.Range("L126").Value = 12174
$ws.Range("M126").Value = -106408.004
$ws.Range("N126").Value = -17114
$ws.Range("H132").Value = 5228.2856
$ws.Range("I132").Value = 3819.6
$ws.Range("J132").Value = 8750
$ws.Range("K132").Value = 11458.8
$ws.Range("L132").Value = 26250
$ws.Range("M132").Value = -8928.799999999999
$ws.Range("N132").Value = -31310

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H81").Value = 6196.933
$ws.Range("I81").Value = 5872.125
$ws.Range("J81").Value = 6568.143
$ws.Range("K81").Value = 11744.25
$ws.Range("L81").Value = 13136.286
$ws.Range("M81").Value = -10683.25
$ws.Range("N81").Value = -15258.286
$ws.Range("H84").Value = 6196.933
$ws.Range("I84").Value = 5872.125
$ws.Range("J84").Value = 6568.143
$ws.Range("K84").Value = 58721.25
$ws.Range("L84").Value = 65681.42999999999
$ws.Range("M84").Value = -53417.25
$ws.Range("N84").Value = -76289.42999999999
$ws.Range("H113").Value = 1077.6842
$ws.Range("I113").Value = 1167.1818
$ws.Range("J113").Value = 954.625
$ws.Range("K113").Value = 3501.5454
$ws.Range("L113").Value = 2863.875
$ws.Range("M113").Value = -1331.5454
$ws.Range("N113").Value = -7203.875
$ws.Range("H121").Value = 48307.668
$ws.Range("J121").Value = 48307.668
$ws.Range("L121").Value = 48307.668
$ws.Range("N121").Value = -51801.668
